# Adding updated GymWorkout data 11/03/2018
#
# Appends 10 new workout-log rows (rows 620-629) to the WeightTraining
# sheet, all recorded on 11/03/2018 (a Sunday) for DateId 71.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: ExerciseId, DateId, ExerciseDate(serial), ExerciseWeek,
#             ExerciseMonth, ExerciseYear, ExerciseDay, ExerciseName,
#             Weight, Sets, Reps, TrainingArea
$newRows = @(
    @(619, 71, 43170, 12, "March", 2018, "Sunday", "Leg Extension",       88,  4, 12, "Legs"),
    @(620, 71, 43170, 12, "March", 2018, "Sunday", "Pec Fly",             100, 4, 8,  "Chest"),
    @(621, 71, 43170, 12, "March", 2018, "Sunday", "Hip abduction",       60,  3, 12, "Legs"),
    @(622, 71, 43170, 12, "March", 2018, "Sunday", "Hip adduction",       65,  3, 12, "Legs"),
    @(623, 71, 43170, 12, "March", 2018, "Sunday", "Seated Row",          90,  4, 8,  "Back"),
    @(624, 71, 43170, 12, "March", 2018, "Sunday", "Heel-taps",           0,   4, 10, "Core"),
    @(625, 71, 43170, 12, "March", 2018, "Sunday", "Raised leg circles",  0,   4, 10, "Core"),
    @(626, 71, 43170, 12, "March", 2018, "Sunday", "Scissors",            0,   4, 12, "Core"),
    @(627, 71, 43170, 12, "March", 2018, "Sunday", "Knee-Pull ins",       0,   4, 10, "Core"),
    @(628, 71, 43170, 12, "March", 2018, "Sunday", "Flitter Kicks",       0,   4, 20, "Core")
)

$startRow = 620
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $data = $newRows[$i]

    $ws.Cells.Item($r, 1).Value  = $data[0]   # ExerciseId
    $ws.Cells.Item($r, 2).Value  = $data[1]   # DateId
    $ws.Cells.Item($r, 3).Value  = $data[2]   # Exercise Date
    $ws.Cells.Item($r, 4).Value  = $data[3]   # Exercise Week
    $ws.Cells.Item($r, 5).Value  = $data[4]   # Exercise Month
    $ws.Cells.Item($r, 6).Value  = $data[5]   # Exercise Year
    $ws.Cells.Item($r, 7).Value  = $data[6]   # Exercise Day
    $ws.Cells.Item($r, 8).Value  = $data[7]   # Exercise Name
    $ws.Cells.Item($r, 9).Value  = $data[8]   # Weight
    $ws.Cells.Item($r, 10).Value = $data[9]   # Sets
    $ws.Cells.Item($r, 11).Value = $data[10]  # Reps
    $ws.Cells.Item($r, 12).Value = $data[11]  # TrainingArea
}

$lastRow = $startRow + $newRows.Count - 1

# Bring the newly added rows into view and select them, mirroring the
# updated scroll position / selection left behind after data entry.
$win = $wb.Windows.Item(1)
$win.ScrollRow = 601
$win.ScrollColumn = 1
$ws.Range("A611:A$lastRow").Select() | Out-Null
